# --- Update auto date fields (datetimeFigureOut) on slide master and all layouts ---
$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$mshapes = $m.Shapes
for ($i = 1; $i -le $mshapes.Count; $i++) {
    $sh = $mshapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "7/3/18"
    }
}

$cls = $m.CustomLayouts
for ($i = 1; $i -le $cls.Count; $i++) {
    $cl = $cls.Item($i)
    $shapes = $cl.Shapes
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "7/3/18"
        }
    }
}

# --- Slide 10 ("General pattern"): rewrite body text to add citation paragraph ---
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(2)
$tr10 = $sh10.TextFrame.TextRange
$tr10.Text = "The technique on previous slide can be generalized to an entire network of multiple devices.`rSee “Abstractions for Network Update”, M. Reitblatt, N. Foster, J. Rexford, C. Schlesinger, D. Walker, SIGCOMM 2012`rIt is a general ‘atomic pointer flipping’ technique that can be applied to smaller subsets of tables.`rThe ‘pointer’ being flipped in previous slide is the blue/red bit, ‘pointing’ at one of two sets of tables.`rFleshed out example starts on next slide."

# Set indent level for citation paragraph (2) and the previously-2nd paragraph (now 4), both lvl=1
$para2 = $tr10.Characters(94, 115)
$para2.IndentLevel = 2
$para4 = $tr10.Characters(312, 107)
$para4.IndentLevel = 2

# Split citation paragraph into 3 runs (so "Reitblatt" is its own run)
$r2a = $tr10.Characters(94, 42)
$r2a.Text = "See “Abstractions for Network Update”, M. "
$r2b = $tr10.Characters(136, 9)
$r2b.Text = "Reitblatt"
$r2c = $tr10.Characters(145, 64)
$r2c.Text = ", N. Foster, J. Rexford, C. Schlesinger, D. Walker, SIGCOMM 2012"

# Split third paragraph into 2 runs ("It is " + rest)
$r3a = $tr10.Characters(210, 6)
$r3a.Text = "It is "
$r3b = $tr10.Characters(216, 95)
$r3b.Text = "a general ‘atomic pointer flipping’ technique that can be applied to smaller subsets of tables."
